$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(29, 1).Value = 111277950
$ws.Cells.Item(29, 16).Value = 'Präststranden, Jmt'
$ws.Cells.Item(29, 17).Value = 569874.8142812594
$ws.Cells.Item(29, 18).Value = 6993100.559414167
$ws.Cells.Item(29, 26).Value = '20:26'
$ws.Cells.Item(29, 28).Value = '20:26'
$ws.Cells.Item(30, 1).Value = 111277538
$ws.Cells.Item(30, 17).Value = 569784.6764437903
$ws.Cells.Item(30, 18).Value = 6992856.400962653
$ws.Cells.Item(31, 1).Value = 111277389
$ws.Cells.Item(31, 17).Value = 569750.3053765292
$ws.Cells.Item(31, 18).Value = 6992912.817455334
$ws.Cells.Item(32, 1).Value = 111278217
$ws.Cells.Item(32, 2).Value = 89686
$ws.Cells.Item(32, 4).Value = 'NT'
$ws.Cells.Item(32, 5).Value = 658
$ws.Cells.Item(32, 6).Value = 'Rosenticka'
$ws.Cells.Item(32, 7).Value = 'Rhodofomes roseus'
$ws.Cells.Item(32, 8).Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(32, 16).Value = 'Singsån, Jmt'
$ws.Cells.Item(32, 17).Value = 569671.7019483433
$ws.Cells.Item(32, 18).Value = 6993040.858867787
$ws.Cells.Item(33, 1).Value = 111277633
$ws.Cells.Item(33, 17).Value = 569802.0407188418
$ws.Cells.Item(33, 18).Value = 6992830.464391444
$ws.Cells.Item(34, 1).Value = 111277392
$ws.Cells.Item(34, 2).Value = 89845
$ws.Cells.Item(34, 4).Value = 'VU'
$ws.Cells.Item(34, 5).Value = 1209
$ws.Cells.Item(34, 6).Value = 'Rynkskinn'
$ws.Cells.Item(34, 7).Value = 'Phlebia centrifuga'
$ws.Cells.Item(34, 8).Value = 'P.Karst.'
$ws.Cells.Item(34, 17).Value = 569750.3053765292
$ws.Cells.Item(34, 18).Value = 6992912.817455334
$ws.Cells.Item(34, 26).Value = '00:00'
$ws.Cells.Item(34, 28).Value = '00:00'
$ws.Cells.Item(35, 1).Value = 111278872
$ws.Cells.Item(35, 2).Value = 78578
$ws.Cells.Item(35, 4).Value = 'NT'
$ws.Cells.Item(35, 5).Value = 6458
$ws.Cells.Item(35, 6).Value = 'Lunglav'
$ws.Cells.Item(35, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(35, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(35, 17).Value = 569296.7869269754
$ws.Cells.Item(35, 18).Value = 6992794.243538878
$ws.Cells.Item(35, 26).Value = '20:26'
$ws.Cells.Item(35, 28).Value = '20:26'
$ws.Cells.Item(37, 1).Value = 111277903
$ws.Cells.Item(37, 2).Value = 96348
$ws.Cells.Item(37, 4).Value = 'VU'
$ws.Cells.Item(37, 5).Value = 220787
$ws.Cells.Item(37, 6).Value = 'Knärot'
$ws.Cells.Item(37, 7).Value = 'Goodyera repens'
$ws.Cells.Item(37, 8).Value = '(L.) R. Br.'
$ws.Cells.Item(37, 16).Value = 'Präststranden, Jmt'
$ws.Cells.Item(37, 17).Value = 569897.0842333297
$ws.Cells.Item(37, 18).Value = 6993078.813114846
$ws.Cells.Item(37, 26).Value = '20:24'
$ws.Cells.Item(37, 28).Value = '20:24'
$ws.Cells.Item(38, 1).Value = 111279094
$ws.Cells.Item(38, 2).Value = 89416
$ws.Cells.Item(38, 4).Value = 'LC'
$ws.Cells.Item(38, 5).Value = 1205
$ws.Cells.Item(38, 6).Value = 'Stor aspticka'
$ws.Cells.Item(38, 7).Value = 'Phellinus populicola'
$ws.Cells.Item(38, 8).Value = 'Niemelä'
$ws.Cells.Item(38, 17).Value = 569279.6199819668
$ws.Cells.Item(38, 18).Value = 6992811.114809629
$ws.Cells.Item(39, 1).Value = 111278492
$ws.Cells.Item(39, 16).Value = 'Ragunda, Jmt'
$ws.Cells.Item(39, 17).Value = 569641.4769454591
$ws.Cells.Item(39, 18).Value = 6992967.635971196
$ws.Cells.Item(40, 1).Value = 111277552
$ws.Cells.Item(40, 2).Value = 89686
$ws.Cells.Item(40, 4).Value = 'NT'
$ws.Cells.Item(40, 5).Value = 658
$ws.Cells.Item(40, 6).Value = 'Rosenticka'
$ws.Cells.Item(40, 7).Value = 'Rhodofomes roseus'
$ws.Cells.Item(40, 8).Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Cells.Item(40, 16).Value = 'Ragunda, Jmt'
$ws.Cells.Item(40, 17).Value = 569770.841244747
$ws.Cells.Item(40, 18).Value = 6992866.083226931
$ws.Cells.Item(40, 26).Value = '00:00'
$ws.Cells.Item(40, 28).Value = '00:00'
$ws.Cells.Item(41, 1).Value = 111279409
$ws.Cells.Item(41, 16).Value = 'Ragunda, Jmt'
$ws.Cells.Item(41, 17).Value = 569443.239979364
$ws.Cells.Item(41, 18).Value = 6992913.042043422
$ws.Cells.Item(41, 26).Value = '21:39'
$ws.Cells.Item(41, 28).Value = '21:39'
